$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.188.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.499.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.82%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.516.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0975"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.57%  "

$ws.Range("E13").Value = "  +1.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.930.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.130.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.60%  "

$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.514.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.411"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.67%  "

$ws.Range("E26").Value = "  +6.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.617.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0790"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.58%  "

$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("E34").Value = "  +4.77%  "

$ws.Range("E35").Value = "  +3.21%  "

$ws.Range("E36").Value = "  +7.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.869"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.619"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.59%  "

$ws.Range("E42").Value = "  +4.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.993"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "260.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +17.20%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0229"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.58%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0914"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.76%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.922.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.74%  "
